# Auto-generated script to update market/profit data cells per commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 197.83333
$ws.Cells.Item(2, 9).Value = 175.46153
$ws.Cells.Item(2, 10).Value = 256
$ws.Cells.Item(2, 11).Value = 175.46153
$ws.Cells.Item(2, 12).Value = 256
$ws.Cells.Item(2, 13).Value = -62.46153000000001
$ws.Cells.Item(2, 14).Value = -482
$ws.Cells.Item(70, 8).Value = 1942.1177
$ws.Cells.Item(70, 9).Value = 1078.1538
$ws.Cells.Item(70, 10).Value = 4750
$ws.Cells.Item(70, 11).Value = 3234.4614
$ws.Cells.Item(70, 12).Value = 14250
$ws.Cells.Item(70, 13).Value = -2964.4614
$ws.Cells.Item(70, 14).Value = -14790
$ws.Cells.Item(73, 8).Value = 1942.1177
$ws.Cells.Item(73, 9).Value = 1078.1538
$ws.Cells.Item(73, 10).Value = 4750
$ws.Cells.Item(73, 11).Value = 3234.4614
$ws.Cells.Item(73, 12).Value = 14250
$ws.Cells.Item(73, 13).Value = -2298.4614
$ws.Cells.Item(73, 14).Value = -16122
$ws.Cells.Item(76, 8).Value = 3055.0571
$ws.Cells.Item(76, 9).Value = 3052.3333
$ws.Cells.Item(76, 10).Value = 3100
$ws.Cells.Item(76, 11).Value = 3052.3333
$ws.Cells.Item(76, 12).Value = 3100
$ws.Cells.Item(76, 13).Value = -2737.3333
$ws.Cells.Item(76, 14).Value = -3730
$ws.Cells.Item(79, 8).Value = 3055.0571
$ws.Cells.Item(79, 9).Value = 3052.3333
$ws.Cells.Item(79, 10).Value = 3100
$ws.Cells.Item(79, 11).Value = 3052.3333
$ws.Cells.Item(79, 12).Value = 3100
$ws.Cells.Item(79, 13).Value = -1960.3333
$ws.Cells.Item(79, 14).Value = -5284
$ws.Cells.Item(80, 8).Value = 234.63158
$ws.Cells.Item(80, 9).Value = 175.35715
$ws.Cells.Item(80, 10).Value = 400.6
$ws.Cells.Item(80, 11).Value = 526.0714499999999
$ws.Cells.Item(80, 12).Value = 1201.8
$ws.Cells.Item(80, 13).Value = 471.9285500000001
$ws.Cells.Item(80, 14).Value = -3197.8
$ws.Cells.Item(83, 8).Value = 234.63158
$ws.Cells.Item(83, 9).Value = 175.35715
$ws.Cells.Item(83, 10).Value = 400.6
$ws.Cells.Item(83, 11).Value = 1578.21435
$ws.Cells.Item(83, 12).Value = 3605.4
$ws.Cells.Item(83, 13).Value = 3413.78565
$ws.Cells.Item(83, 14).Value = -13589.4
$ws.Cells.Item(98, 8).Value = 2485
$ws.Cells.Item(98, 9).Value = 2596.742
$ws.Cells.Item(98, 10).Value = 753
$ws.Cells.Item(98, 11).Value = 2596.742
$ws.Cells.Item(98, 12).Value = 753
$ws.Cells.Item(98, 13).Value = -1098.742
$ws.Cells.Item(98, 14).Value = -3749
$ws.Cells.Item(106, 8).Value = 4252.5
$ws.Cells.Item(106, 9).Value = 4336.6665
$ws.Cells.Item(106, 10).Value = 4000
$ws.Cells.Item(106, 11).Value = 4336.6665
$ws.Cells.Item(106, 12).Value = 4000
$ws.Cells.Item(106, 13).Value = -3705.6665
$ws.Cells.Item(106, 14).Value = -5262
$ws.Cells.Item(107, 8).Value = 1063.8334
$ws.Cells.Item(107, 10).Value = 1100.3636
$ws.Cells.Item(107, 12).Value = 1100.3636
$ws.Cells.Item(107, 14).Value = -4940.3636
$ws.Cells.Item(111, 8).Value = 1000
$ws.Cells.Item(111, 9).Value = 1000
$ws.Cells.Item(111, 10).Value = 0
$ws.Cells.Item(111, 11).Value = 3000
$ws.Cells.Item(111, 12).Value = 0
$ws.Cells.Item(111, 13).Value = 67
$ws.Cells.Item(111, 14).ClearContents()
$ws.Cells.Item(116, 8).Value = 3091.4707
$ws.Cells.Item(116, 9).Value = 1980.5
$ws.Cells.Item(116, 10).Value = 4678.5713
$ws.Cells.Item(116, 11).Value = 1980.5
$ws.Cells.Item(116, 12).Value = 4678.5713
$ws.Cells.Item(116, 13).Value = 1461.5
$ws.Cells.Item(116, 14).Value = -11562.5713
$ws.Cells.Item(122, 8).Value = 2485
$ws.Cells.Item(122, 9).Value = 2596.742
$ws.Cells.Item(122, 10).Value = 753
$ws.Cells.Item(122, 11).Value = 7790.226000000001
$ws.Cells.Item(122, 12).Value = 2259
$ws.Cells.Item(122, 13).Value = -5340.226000000001
$ws.Cells.Item(122, 14).Value = -7159
$ws.Cells.Item(132, 8).Value = 229797.45
$ws.Cells.Item(132, 9).Value = 252692.2
$ws.Cells.Item(132, 10).Value = 850
$ws.Cells.Item(132, 11).Value = 758076.6000000001
$ws.Cells.Item(132, 12).Value = 2550
$ws.Cells.Item(132, 13).Value = -755546.6000000001
$ws.Cells.Item(132, 14).Value = -7610

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 3012.125
$ws.Cells.Item(88, 9).Value = 3048.3333
$ws.Cells.Item(88, 10).Value = 2903.5
$ws.Cells.Item(88, 11).Value = 3048.3333
$ws.Cells.Item(88, 12).Value = 2903.5
$ws.Cells.Item(88, 13).Value = -2642.3333
$ws.Cells.Item(88, 14).Value = -3715.5
$ws.Cells.Item(91, 8).Value = 3012.125
$ws.Cells.Item(91, 9).Value = 3048.3333
$ws.Cells.Item(91, 10).Value = 2903.5
$ws.Cells.Item(91, 11).Value = 3048.3333
$ws.Cells.Item(91, 12).Value = 2903.5
$ws.Cells.Item(91, 13).Value = -1644.3333
$ws.Cells.Item(91, 14).Value = -5711.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2265.0435
$ws.Cells.Item(86, 9).Value = 1820.4
$ws.Cells.Item(86, 10).Value = 3098.75
$ws.Cells.Item(86, 11).Value = 1820.4
$ws.Cells.Item(86, 12).Value = 3098.75
$ws.Cells.Item(86, 13).Value = -697.4000000000001
$ws.Cells.Item(86, 14).Value = -5344.75
$ws.Cells.Item(89, 8).Value = 2265.0435
$ws.Cells.Item(89, 9).Value = 1820.4
$ws.Cells.Item(89, 10).Value = 3098.75
$ws.Cells.Item(89, 11).Value = 9102
$ws.Cells.Item(89, 12).Value = 15493.75
$ws.Cells.Item(89, 13).Value = -3486
$ws.Cells.Item(89, 14).Value = -26725.75
$ws.Cells.Item(134, 8).Value = 41765.76
$ws.Cells.Item(134, 9).Value = 41765.76
$ws.Cells.Item(134, 11).Value = 125297.28
$ws.Cells.Item(134, 13).Value = -122762.28

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1543.4667
$ws.Cells.Item(31, 9).Value = 1563.0834
$ws.Cells.Item(31, 10).Value = 1465
$ws.Cells.Item(31, 11).Value = 1563.0834
$ws.Cells.Item(31, 12).Value = 1465
$ws.Cells.Item(31, 13).Value = -1268.0834
$ws.Cells.Item(31, 14).Value = -2055
$ws.Cells.Item(34, 8).Value = 1543.4667
$ws.Cells.Item(34, 9).Value = 1563.0834
$ws.Cells.Item(34, 10).Value = 1465
$ws.Cells.Item(34, 11).Value = 1563.0834
$ws.Cells.Item(34, 12).Value = 1465
$ws.Cells.Item(34, 13).Value = -1361.0834
$ws.Cells.Item(34, 14).Value = -1869
$ws.Cells.Item(107, 8).Value = 551.4737
$ws.Cells.Item(107, 9).Value = 700.1818
$ws.Cells.Item(107, 10).Value = 347
$ws.Cells.Item(107, 11).Value = 700.1818
$ws.Cells.Item(107, 12).Value = 347
$ws.Cells.Item(107, 13).Value = 1219.8182
$ws.Cells.Item(107, 14).Value = -4187
$ws.Cells.Item(118, 8).Value = 49741.11
$ws.Cells.Item(118, 10).Value = 49741.11
$ws.Cells.Item(118, 12).Value = 49741.11
$ws.Cells.Item(118, 14).Value = -53055.11

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(13, 8).Value = 129.5
$ws.Cells.Item(13, 9).Value = 129.5
$ws.Cells.Item(13, 11).Value = 388.5
$ws.Cells.Item(13, 13).Value = -220.5
$ws.Cells.Item(35, 8).Value = 1110
$ws.Cells.Item(35, 10).Value = 1110
$ws.Cells.Item(35, 12).Value = 3330
$ws.Cells.Item(35, 14).Value = -3906
$ws.Cells.Item(57, 8).Value = 3000
$ws.Cells.Item(57, 9).Value = 0
$ws.Cells.Item(57, 10).Value = 3000
$ws.Cells.Item(57, 11).Value = 0
$ws.Cells.Item(57, 12).Value = 9000
$ws.Cells.Item(57, 14).Value = -10118
$ws.Cells.Item(57, 13).ClearContents()
$ws.Cells.Item(93, 8).Value = 17679.8
$ws.Cells.Item(93, 10).Value = 2099.75
$ws.Cells.Item(93, 12).Value = 6299.25
$ws.Cells.Item(93, 14).Value = -10043.25
$ws.Cells.Item(94, 8).Value = 2308
$ws.Cells.Item(94, 10).Value = 2800
$ws.Cells.Item(94, 12).Value = 8400
$ws.Cells.Item(94, 14).Value = -9752
$ws.Cells.Item(95, 8).Value = 3300
$ws.Cells.Item(95, 10).Value = 3300
$ws.Cells.Item(95, 12).Value = 9900
$ws.Cells.Item(95, 14).Value = -14018
$ws.Cells.Item(99, 8).Value = 2027
$ws.Cells.Item(99, 10).Value = 3028
$ws.Cells.Item(99, 12).Value = 9084
$ws.Cells.Item(99, 14).Value = -13576
$ws.Cells.Item(101, 8).Value = 5401
$ws.Cells.Item(101, 10).Value = 6792.6665
$ws.Cells.Item(101, 12).Value = 20377.9995
$ws.Cells.Item(101, 14).Value = -25245.9995
$ws.Cells.Item(102, 8).Value = 9007.25
$ws.Cells.Item(102, 10).Value = 9007.25
$ws.Cells.Item(102, 12).Value = 27021.75
$ws.Cells.Item(102, 14).Value = -31889.75
$ws.Cells.Item(106, 8).Value = 5480
$ws.Cells.Item(106, 10).Value = 5480
$ws.Cells.Item(106, 12).Value = 16440
$ws.Cells.Item(106, 14).Value = -18332
$ws.Cells.Item(117, 8).Value = 709.6667
$ws.Cells.Item(117, 9).Value = 709.6667
$ws.Cells.Item(117, 11).Value = 2129.0001
$ws.Cells.Item(117, 13).Value = 1312.9999
$ws.Cells.Item(118, 8).Value = 1357.8
$ws.Cells.Item(118, 9).Value = 709.75
$ws.Cells.Item(118, 11).Value = 2129.25
$ws.Cells.Item(118, 13).Value = -886.25
$ws.Cells.Item(131, 8).Value = 3396.75
$ws.Cells.Item(131, 10).Value = 2413.4285
$ws.Cells.Item(131, 12).Value = 7240.2855
$ws.Cells.Item(131, 14).Value = -17320.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 1659.3684
$ws.Cells.Item(132, 9).Value = 1253.9231
$ws.Cells.Item(132, 10).Value = 2537.8333
$ws.Cells.Item(132, 11).Value = 3761.7693
$ws.Cells.Item(132, 12).Value = 7613.499899999999
$ws.Cells.Item(132, 13).Value = -1231.7693
$ws.Cells.Item(132, 14).Value = -12673.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1978.0952
$ws.Cells.Item(16, 9).Value = 2239.4443
$ws.Cells.Item(16, 10).Value = 410
$ws.Cells.Item(16, 11).Value = 2239.4443
$ws.Cells.Item(16, 12).Value = 410
$ws.Cells.Item(16, 13).Value = -2069.4443
$ws.Cells.Item(16, 14).Value = -750
$ws.Cells.Item(136, 8).Value = 1799.1852
$ws.Cells.Item(136, 9).Value = 1071.8334
$ws.Cells.Item(136, 10).Value = 3253.889
$ws.Cells.Item(136, 11).Value = 3215.5002
$ws.Cells.Item(136, 12).Value = 9761.667000000001
$ws.Cells.Item(136, 13).Value = -665.5001999999999
$ws.Cells.Item(136, 14).Value = -14861.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 14462111
$ws.Cells.Item(2, 10).Value = 15737000
$ws.Cells.Item(2, 12).Value = 15737000
$ws.Cells.Item(2, 14).Value = -15737224
$ws.Cells.Item(132, 8).Value = 3363.8965
$ws.Cells.Item(132, 9).Value = 3817.05
$ws.Cells.Item(132, 10).Value = 2356.889
$ws.Cells.Item(132, 11).Value = 11451.15
$ws.Cells.Item(132, 12).Value = 7070.667
$ws.Cells.Item(132, 13).Value = -8921.150000000001
$ws.Cells.Item(132, 14).Value = -12130.667
